$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new forecast-origin row (year 2007) is inserted at row 2, which
# pushes the existing rows down by one (old row 2 -> new row 3, etc.).
$ws.Rows.Item(2).Insert()

# Copy the number formatting/borders from the row below (the old row 2,
# now shifted to row 3) onto the freshly inserted row 2, so the new row
# matches the look of the rest of the table (date style in column A).
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-simulated / bugfixed naive forecaster values for every row.
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 11.13090654781819
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 9.396507498425466

$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 4.672550446571067
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 4.422525088127283

$ws.Cells.Item(4, 1).Value = 40130
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = -14.45332333832743
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).Value = -2.928447329610073

$ws.Cells.Item(5, 1).Value = 40494
$ws.Cells.Item(5, 2).Value = 2010
$ws.Cells.Item(5, 3).Value = 8.600536527919633
$ws.Cells.Item(5, 4).Value = 2011
$ws.Cells.Item(5, 5).Value = 6.303897256856628

$ws.Cells.Item(6, 1).Value = 40862
$ws.Cells.Item(6, 2).Value = 2011
$ws.Cells.Item(6, 3).Value = 10.25770250047622
$ws.Cells.Item(6, 4).Value = 2012
$ws.Cells.Item(6, 5).Value = 10.22374275635105

$ws.Cells.Item(7, 1).Value = 41228
$ws.Cells.Item(7, 2).Value = 2012
$ws.Cells.Item(7, 3).Value = 4.639893381363169
$ws.Cells.Item(7, 4).Value = 2013
$ws.Cells.Item(7, 5).Value = 8.174613408931286

$ws.Cells.Item(8, 1).Value = 41592
$ws.Cells.Item(8, 2).Value = 2013
$ws.Cells.Item(8, 3).Value = 0.3058963467304165
$ws.Cells.Item(8, 4).Value = 2014
$ws.Cells.Item(8, 5).Value = 2.429116709932622

$ws.Cells.Item(9, 1).Value = 41957
$ws.Cells.Item(9, 2).Value = 2014
$ws.Cells.Item(9, 3).Value = 4.068173739091874
$ws.Cells.Item(9, 4).Value = 2015
$ws.Cells.Item(9, 5).Value = 3.9413000500929

$ws.Cells.Item(10, 1).Value = 42321
$ws.Cells.Item(10, 2).Value = 2015
$ws.Cells.Item(10, 3).Value = 4.984288257750213
$ws.Cells.Item(10, 4).Value = 2016
$ws.Cells.Item(10, 5).Value = 4.188839638544284

$ws.Cells.Item(11, 1).Value = 42689
$ws.Cells.Item(11, 2).Value = 2016
$ws.Cells.Item(11, 3).Value = 1.878184267712912
$ws.Cells.Item(11, 4).Value = 2017
$ws.Cells.Item(11, 5).Value = 2.514670279852349

$ws.Cells.Item(12, 1).Value = 43053
$ws.Cells.Item(12, 2).Value = 2017
$ws.Cells.Item(12, 3).Value = 4.695933104194339
$ws.Cells.Item(12, 4).Value = 2018
$ws.Cells.Item(12, 5).Value = 4.5579527192392

$ws.Cells.Item(13, 1).Value = 43418
$ws.Cells.Item(13, 2).Value = 2018
$ws.Cells.Item(13, 3).Value = 4.892602738886098
$ws.Cells.Item(13, 4).Value = 2019
$ws.Cells.Item(13, 5).Value = 1.957202207503861

$ws.Cells.Item(14, 1).Value = 43783
$ws.Cells.Item(14, 2).Value = 2019
$ws.Cells.Item(14, 3).Value = 0.8049382522247184
$ws.Cells.Item(14, 4).Value = 2020
$ws.Cells.Item(14, 5).Value = 2.267257846564918

$ws.Cells.Item(15, 1).Value = 44159
$ws.Cells.Item(15, 2).Value = 2020
$ws.Cells.Item(15, 3).Value = -8.78417389973717
$ws.Cells.Item(15, 4).Value = 2021
$ws.Cells.Item(15, 5).Value = 2.199380357735481

$ws.Cells.Item(16, 1).Value = 44525
$ws.Cells.Item(16, 2).Value = 2021
$ws.Cells.Item(16, 3).Value = 5.110501195359984
$ws.Cells.Item(16, 4).Value = 2022
$ws.Cells.Item(16, 5).Value = 0.3515918738370427

$ws.Cells.Item(17, 1).Value = 44890
$ws.Cells.Item(17, 2).Value = 2022
$ws.Cells.Item(17, 3).Value = 5.120680133083599
$ws.Cells.Item(17, 4).Value = 2023
$ws.Cells.Item(17, 5).Value = 5.934275247805543

$ws.Cells.Item(18, 1).Value = 45254
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = -0.5532735011319234
$ws.Cells.Item(18, 4).Value = 2024
$ws.Cells.Item(18, 5).Value = -1.846917864698006

$ws.Cells.Item(19, 1).Value = 45618
$ws.Cells.Item(19, 2).Value = 2024
$ws.Cells.Item(19, 3).Value = -1.069674659641462
$ws.Cells.Item(19, 4).Value = 2025
$ws.Cells.Item(19, 5).Value = -0.7986414110784379

